$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14.00156945261478
$ws.Range("D2").Value = 3.411382621167844
$ws.Range("E2").Value = 8.985785499056439
$ws.Range("F2").Value = 57.97724184698689
$ws.Range("G2").Value = 3.804184768701249
$ws.Range("I2").Value = 39.59812849506417
$ws.Range("J2").Value = 9.497483456414422
$ws.Range("K2").Value = 17.23823850884696
$ws.Range("L2").Value = 13.33790825060601
$ws.Range("M2").Value = 17.275477066505
$ws.Range("B3").Value = 14.00841760818214
$ws.Range("D3").Value = 3.35810767224216
$ws.Range("E3").Value = 8.938955472372626
$ws.Range("F3").Value = 57.75394045640792
$ws.Range("G3").Value = 3.807766946923716
$ws.Range("I3").Value = 39.51318809476727
$ws.Range("J3").Value = 9.486174597047603
$ws.Range("K3").Value = 17.15940188278502
$ws.Range("L3").Value = 13.34742435539122
$ws.Range("M3").Value = 17.30406441821459
$ws.Range("B4").Value = 14.01634070560273
$ws.Range("D4").Value = 3.32522786256811
$ws.Range("E4").Value = 8.909455892169236
$ws.Range("F4").Value = 57.62428570162286
$ws.Range("G4").Value = 3.810081304246659
$ws.Range("I4").Value = 39.4643115185903
$ws.Range("J4").Value = 9.478962381337086
$ws.Range("K4").Value = 17.11597361532462
$ws.Range("L4").Value = 13.35568441268197
$ws.Range("M4").Value = 17.32473920448031
$ws.Range("B5").Value = 14.02050562493183
$ws.Range("D5").Value = 3.311797928676363
$ws.Range("E5").Value = 8.897247874322593
$ws.Range("F5").Value = 57.57334830261632
$ws.Range("G5").Value = 3.811053418081809
$ws.Range("I5").Value = 39.44522298514869
$ws.Range("J5").Value = 9.475955493300429
$ws.Range("K5").Value = 17.09954301640233
$ws.Range("L5").Value = 13.35965828275586
$ws.Range("M5").Value = 17.33394951994304
$ws.Range("B6").Value = 14.02125377533464
$ws.Range("D6").Value = 3.309566348654529
$ws.Range("E6").Value = 8.895209490036162
$ws.Range("F6").Value = 57.56500547259665
$ws.Range("G6").Value = 3.811216591138206
$ws.Range("I6").Value = 39.44210355756469
$ws.Range("J6").Value = 9.475452094176996
$ws.Range("K6").Value = 17.09689164096553
$ws.Range("L6").Value = 13.36035484934633
$ws.Range("M6").Value = 17.33552630889949
$ws.Range("B7").Value = 14.0163930833698
$ws.Range("D7").Value = 3.325046852924037
$ws.Range("E7").Value = 8.909292005388677
$ws.Range("F7").Value = 57.62359102757863
$ws.Range("G7").Value = 3.810094296977782
$ws.Range("I7").Value = 39.46405072100477
$ws.Range("J7").Value = 9.478922104541764
$ws.Range("K7").Value = 17.11574687895408
$ws.Range("L7").Value = 13.35573554473451
$ws.Range("M7").Value = 17.32486023886053
$ws.Range("B8").Value = 14.0031598294818
$ws.Range("D8").Value = 3.393055196541868
$ws.Range("E8").Value = 8.969793778120431
$ws.Range("F8").Value = 57.89871549553014
$ws.Range("G8").Value = 3.805396122944154
$ws.Range("I8").Value = 39.56816400801337
$ws.Range("J8").Value = 9.493639696845824
$ws.Range("K8").Value = 17.21003213398713
$ws.Range("L8").Value = 13.34068779172387
$ws.Range("M8").Value = 17.28468610951217
$ws.Range("B9").Value = 14.00663884332556
$ws.Range("D9").Value = 3.537067579575966
$ws.Range("E9").Value = 9.082486924485115
$ws.Range("F9").Value = 58.49624373529269
$ws.Range("G9").Value = 3.797089768784492
$ws.Range("I9").Value = 39.79809830220122
$ws.Range("J9").Value = 9.520385418572719
$ws.Range("K9").Value = 17.43367000419589
$ws.Range("L9").Value = 13.33035134665167
$ws.Range("M9").Value = 17.23066820266022
$ws.Range("B10").Value = 14.02700560832056
$ws.Range("D10").Value = 3.669201935085049
$ws.Range("E10").Value = 9.161614582433234
$ws.Range("F10").Value = 58.96898703541361
$ws.Range("G10").Value = 3.791533136617857
$ws.Range("I10").Value = 39.98240463735124
$ws.Range("J10").Value = 9.538767631434498
$ws.Range("K10").Value = 17.62046863615712
$ws.Range("L10").Value = 13.33442746642768
$ws.Range("M10").Value = 17.20606061896505
$ws.Range("B11").Value = 14.0400995364168
$ws.Range("D11").Value = 3.727582316256226
$ws.Range("E11").Value = 9.196806897870131
$ws.Range("F11").Value = 59.19099149439172
$ws.Range("G11").Value = 3.789122413032257
$ws.Range("I11").Value = 40.06950654388118
$ws.Range("J11").Value = 9.546859565417707
$ws.Range("K11").Value = 17.71005755304123
$ws.Range("L11").Value = 13.33880822979986
$ws.Range("M11").Value = 17.19813380718543
$ws.Range("B12").Value = 14.04560428006895
$ws.Range("D12").Value = 3.749429080753515
$ws.Range("E12").Value = 9.210016956241097
$ws.Range("F12").Value = 59.27602275347849
$ws.Range("G12").Value = 3.788226251485608
$ws.Range("I12").Value = 40.10295046089338
$ws.Range("J12").Value = 9.549885268288543
$ws.Range("K12").Value = 17.74462063883683
$ws.Range("L12").Value = 13.34082937767179
$ws.Range("M12").Value = 17.19560111506263
$ws.Range("B13").Value = 14.04439450907713
$ws.Range("D13").Value = 3.744735762378939
$ws.Range("E13").Value = 9.207177131451035
$ws.Range("F13").Value = 59.25766752057375
$ws.Range("G13").Value = 3.788418513545088
$ws.Range("I13").Value = 40.0957273827381
$ws.Range("J13").Value = 9.549235337865193
$ws.Range("K13").Value = 17.73714891977768
$ws.Range("L13").Value = 13.34037799540276
$ws.Range("M13").Value = 17.19612573058849
$ws.Range("B14").Value = 14.04054149144631
$ws.Range("D14").Value = 3.729384951698306
$ws.Range("E14").Value = 9.197896038035054
$ws.Range("F14").Value = 59.1979679827293
$ws.Range("G14").Value = 3.789048350615084
$ws.Range("I14").Value = 40.07224882004214
$ws.Range("J14").Value = 9.547109262309805
$ws.Range("K14").Value = 17.71288844982997
$ws.Range("L14").Value = 13.33896725770209
$ws.Range("M14").Value = 17.19791604648713
$ws.Range("B15").Value = 14.03825242335951
$ws.Range("D15").Value = 3.71994785442909
$ws.Range("E15").Value = 9.192195897918468
$ws.Range("F15").Value = 59.16152458855978
$ws.Range("G15").Value = 3.789436319426245
$ws.Range("I15").Value = 40.0579272087395
$ws.Range("J15").Value = 9.545801970538417
$ws.Range("K15").Value = 17.69811047625808
$ws.Range("L15").Value = 13.33815028035344
$ws.Range("M15").Value = 17.1990737198456
$ws.Range("B16").Value = 14.02622658786449
$ws.Range("D16").Value = 3.665350806357837
$ws.Range("E16").Value = 9.159298432762924
$ws.Range("F16").Value = 58.95461530711512
$ws.Range("G16").Value = 3.791693029209505
$ws.Range("I16").Value = 39.97677721071386
$ws.Range("J16").Value = 9.538233394339066
$ws.Range("K16").Value = 17.6147044550865
$ws.Range("L16").Value = 13.33419192724387
$ws.Range("M16").Value = 17.20664433536225
$ws.Range("B17").Value = 14.01982697566237
$ws.Range("D17").Value = 3.631405720146651
$ws.Range("E17").Value = 9.138910458663222
$ws.Range("F17").Value = 58.82943801708736
$ws.Range("G17").Value = 3.793107347175072
$ws.Range("I17").Value = 39.92782342563392
$ws.Range("J17").Value = 9.533521335128963
$ws.Range("K17").Value = 17.5647012783075
$ws.Range("L17").Value = 13.33241001864457
$ws.Range("M17").Value = 17.21212496984924
$ws.Range("B18").Value = 14.01650677937953
$ws.Range("D18").Value = 3.61171939734465
$ws.Range("E18").Value = 9.127108289080073
$ws.Range("F18").Value = 58.75809592333423
$ws.Range("G18").Value = 3.793931845512458
$ws.Range("I18").Value = 39.89997389356783
$ws.Range("J18").Value = 9.530785629057009
$ws.Range("K18").Value = 17.5363766211589
$ws.Range("L18").Value = 13.33162299661204
$ws.Range("M18").Value = 17.21558487893893
$ws.Range("B19").Value = 14.01544468739161
$ws.Range("D19").Value = 3.605026507165222
$ws.Range("E19").Value = 9.123099332667671
$ws.Range("F19").Value = 58.73405457094887
$ws.Range("G19").Value = 3.794212902230031
$ws.Range("I19").Value = 39.89059756462963
$ws.Range("J19").Value = 9.529854984409925
$ws.Range("K19").Value = 17.52686199733973
$ws.Range("L19").Value = 13.33139741026177
$ws.Range("M19").Value = 17.21680919619459
$ws.Range("B20").Value = 14.02047092420188
$ws.Range("D20").Value = 3.635036090433114
$ws.Range("E20").Value = 9.141088619229762
$ws.Range("F20").Value = 58.84269566068687
$ws.Range("G20").Value = 3.792955650806741
$ws.Range("I20").Value = 39.93300288578762
$ws.Range("J20").Value = 9.534025573939521
$ws.Range("K20").Value = 17.56997927532188
$ws.Range("L20").Value = 13.33257509278098
$ws.Range("M20").Value = 17.21150971768565
$ws.Range("B21").Value = 14.04165842469835
$ws.Range("D21").Value = 3.733901021552629
$ws.Range("E21").Value = 9.200625290103487
$ws.Range("F21").Value = 59.21547734831929
$ws.Range("G21").Value = 3.788862899021773
$ws.Range("I21").Value = 40.07913262030557
$ws.Range("J21").Value = 9.547734785557598
$ws.Range("K21").Value = 17.71999724157624
$ws.Range("L21").Value = 13.33937180377148
$ws.Range("M21").Value = 17.19737746556788
$ws.Range("B22").Value = 14.05868791454361
$ws.Range("D22").Value = 3.796990966682222
$ws.Range("E22").Value = 9.23885638772005
$ws.Range("F22").Value = 59.46470670105243
$ws.Range("G22").Value = 3.786285501921789
$ws.Range("I22").Value = 40.17731562604527
$ws.Range("J22").Value = 9.556470074346427
$ws.Range("K22").Value = 17.82174694207523
$ws.Range("L22").Value = 13.3459246628186
$ws.Range("M22").Value = 17.19087453091263
$ws.Range("B23").Value = 14.04930932463769
$ws.Range("D23").Value = 3.763461864570935
$ws.Range("E23").Value = 9.218514241136555
$ws.Range("F23").Value = 59.33118885494756
$ws.Range("G23").Value = 3.787652223110474
$ws.Range("I23").Value = 40.12467133241494
$ws.Range("J23").Value = 9.551828319501137
$ws.Range("K23").Value = 17.76711106768512
$ws.Range("L23").Value = 13.34223454442107
$ws.Range("M23").Value = 17.19409548352227
$ws.Range("B24").Value = 14.02017867643161
$ws.Range("D24").Value = 3.633395331264774
$ws.Range("E24").Value = 9.140104124634634
$ws.Range("F24").Value = 58.83669993348403
$ws.Range("G24").Value = 3.793024197254227
$ws.Range("I24").Value = 39.93066033335114
$ws.Range("J24").Value = 9.533797690714939
$ws.Range("K24").Value = 17.56759177385003
$ws.Range("L24").Value = 13.33249972314198
$ws.Range("M24").Value = 17.2117869105662
$ws.Range("B25").Value = 14.00255835681988
$ws.Range("D25").Value = 3.489269735744431
$ws.Range("E25").Value = 9.052641630786223
$ws.Range("F25").Value = 58.32854612258576
$ws.Range("G25").Value = 3.799240481752322
$ws.Range("I25").Value = 39.73316872403299
$ws.Range("J25").Value = 9.5133742138741
$ws.Range("K25").Value = 17.36913186187015
$ws.Range("L25").Value = 13.33109613759463
$ws.Range("M25").Value = 17.24263154806314
